$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new daily data row (row 95) mirroring the existing table layout:
# A = date (kept as text, matching the existing "yyyy/mm/dd" text entries),
# B = weekday (Japanese), C = hour (number), D = ranking (number)

# Force column A to be entered as text so Excel doesn't auto-convert the
# "2025/10/12" string into a date serial value, then restore the default
# "Normal" style so the cell ends up unstyled (same as the rest of the
# data rows, e.g. A94) rather than picking up a new text-format style.
$ws.Cells.Item(95, 1).NumberFormat = "@"
$ws.Cells.Item(95, 1).Value = "2025/10/12"
$ws.Cells.Item(95, 1).Style = "Normal"

$ws.Cells.Item(95, 2).Value = "日"
$ws.Cells.Item(95, 3).Value = 18
$ws.Cells.Item(95, 4).Value = 201
